$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing AgTests (H) / AgPosit (I) values for rows 294-333
$ws.Range("H294").Value = 91360
$ws.Range("I294").Value = 5061

$ws.Range("H300").Value = 71698
$ws.Range("I300").Value = 7031

$ws.Range("H302").Value = 74114
$ws.Range("I302").Value = 5365

$ws.Range("H303").Value = 10423
$ws.Range("I303").Value = 666

$ws.Range("H311").Value = 36069
$ws.Range("I311").Value = 1337

$ws.Range("H312").Value = 40269

$ws.Range("H317").Value = 62528
$ws.Range("I317").Value = 2131

$ws.Range("H319").Value = 56550
$ws.Range("I319").Value = 1787

$ws.Range("H320").Value = 87570
$ws.Range("I320").Value = 3964

$ws.Range("H321").Value = 90401
$ws.Range("I321").Value = 2788

$ws.Range("H322").Value = 106230
$ws.Range("I322").Value = 2280

$ws.Range("H323").Value = 149018
$ws.Range("I323").Value = 2307

$ws.Range("H324").Value = 230440
$ws.Range("I324").Value = 2636

$ws.Range("H325").Value = 698532
$ws.Range("I325").Value = 5756

$ws.Range("H326").Value = 414294
$ws.Range("I326").Value = 3641

$ws.Range("H327").Value = 236946
$ws.Range("I327").Value = 3536

$ws.Range("H328").Value = 180211
$ws.Range("I328").Value = 2634

$ws.Range("H329").Value = 82184

$ws.Range("H330").Value = 70555
$ws.Range("I330").Value = 1970

$ws.Range("H331").Value = 145637
$ws.Range("I331").Value = 2490

$ws.Range("H332").Value = 388839
$ws.Range("I332").Value = 3758

$ws.Range("H333").Value = 238056
$ws.Range("I333").Value = 2533

# Append new row 334 with the latest day's data
$ws.Range("A334").Value = 44228
$ws.Range("B334").Value = 252094
$ws.Range("C334").Value = 222642
$ws.Range("D334").Value = 24668
$ws.Range("E334").Value = 9504
$ws.Range("F334").Value = 1737
$ws.Range("G334").Value = 4784
$ws.Range("H334").Value = 184663
$ws.Range("I334").Value = 3242

# Match the date number format used by the rest of column A
$ws.Range("A334").NumberFormat = $ws.Range("A333").NumberFormat
